$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("listOfCities")

# Update NDTV LiveTemp (column C), Weather Condition NDTV (column F),
# Humidity NDTV (column I), and Wind NDTV (column L) values for rows 2-7.

$ws.Range("C2").Value = "35"
$ws.Range("F2").Value = "Humid and Overcast`n"
$ws.Range("I2").Value = "59"
$ws.Range("L2").Value = "1.23"

$ws.Range("C3").Value = "32"
$ws.Range("F3").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I3").Value = "65"
$ws.Range("L3").Value = "1.60"

$ws.Range("C4").Value = "33"
$ws.Range("F4").Value = "Humid and Overcast`n"
$ws.Range("I4").Value = "60"
$ws.Range("L4").Value = "1.61"

$ws.Range("C5").Value = "31"
$ws.Range("F5").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I5").Value = "62"
$ws.Range("L5").Value = "0.94"

$ws.Range("C6").Value = "30"
$ws.Range("F6").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I6").Value = "77"
$ws.Range("L6").Value = "1.66"

$ws.Range("C7").Value = "34"
$ws.Range("F7").Value = "Humid and Overcast`n"
$ws.Range("I7").Value = "64"
$ws.Range("L7").Value = "2.42"

# Re-autofit rows 2-7 so the newline in column F doesn't leave a stray
# explicit row height behind (matches the original, unset row heights).
$ws.Range("2:7").EntireRow.AutoFit()

# Update the saved selection to C11 (single cell), matching the edited file.
$ws.Range("C11").Select()
